# This script applies the changes described by the commit:
# "Changing t test extractor to be slightly broader, building tests for said same"
#
# It fills in new test-extraction results for rows 68-77 of the
# 'manually_curated_test_set' sheet (replacing the RANDBETWEEN helper
# formulas that used to live in columns L/N for those rows with concrete
# extracted values across columns D-P), and updates the view state
# (zoom level & selection) to match where the author had scrolled to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manually_curated_test_set")

# ---------------------------------------------------------------------
# Row 68
# ---------------------------------------------------------------------
$ws.Range("D68").Value = $false
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = $false
$ws.Range("G68").Value = "NA"
$ws.Range("H68").Value = "NA"
$ws.Range("I68").Value = $false
$ws.Range("J68").Value = "NA"
$ws.Range("K68").Value = 1
$ws.Range("L68").Value = 1
$ws.Range("M68").Value = 4
$ws.Range("N68").Value = 4
$ws.Range("O68").Value = "or in self-administration of saccharin-sweetened water (Figures 1c; t(18)=0.83; P=0.42). Note that responding levels for saccharine are equivalent to those for alcohol in the dependent group."
$ws.Range("P68").Value = "results"

# ---------------------------------------------------------------------
# Row 69
# ---------------------------------------------------------------------
$ws.Range("D69").Value = $false
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = $false
$ws.Range("G69").Value = 80
$ws.Range("H69").Value = "NA"
$ws.Range("I69").Value = $false
$ws.Range("J69").Value = "NA"
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = "NA"
$ws.Range("M69").Value = "NA"
$ws.Range("N69").Value = "NA"
$ws.Range("O69").Value = "NA"
$ws.Range("P69").Value = "NA"

# ---------------------------------------------------------------------
# Row 70
# ---------------------------------------------------------------------
$ws.Range("D70").Value = $false
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = $false
$ws.Range("G70").Value = 200
$ws.Range("H70").Value = "NA"
$ws.Range("I70").Value = $false
$ws.Range("J70").Value = "NA"
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = "NA"
$ws.Range("M70").Value = "NA"
$ws.Range("N70").Value = "NA"
$ws.Range("O70").Value = "NA"
$ws.Range("P70").Value = "NA"

# ---------------------------------------------------------------------
# Row 71
# ---------------------------------------------------------------------
$ws.Range("D71").Value = $false
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = $false
$ws.Range("G71").Value = 41
$ws.Range("H71").Value = $true
$ws.Range("I71").Value = $false
$ws.Range("J71").Value = "NA"
$ws.Range("K71").Value = 5
$ws.Range("L71").Value = 2
$ws.Range("M71").Value = 1
$ws.Range("N71").Value = 1
$ws.Range("O71").Value = "The interaction effect between serum cortisol and WMV of the left MTG was not significant (t = 0.698, p = 0.490) after adjusted for the effect of age, years of education and gender."
$ws.Range("P71").Value = "results"

# ---------------------------------------------------------------------
# Row 72
# ---------------------------------------------------------------------
$ws.Range("D72").Value = $false
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = $false
$ws.Range("G72").Value = "NA"
$ws.Range("H72").Value = "NA"
$ws.Range("I72").Value = $false
$ws.Range("J72").Value = "NA"
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = "NA"
$ws.Range("M72").Value = "NA"
$ws.Range("N72").Value = "NA"
$ws.Range("O72").Value = "NA"
$ws.Range("P72").Value = "NA"

# ---------------------------------------------------------------------
# Row 73
# ---------------------------------------------------------------------
$ws.Range("D73").Value = $false
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = $false
$ws.Range("G73").Value = "NA"
$ws.Range("H73").Value = "NA"
$ws.Range("I73").Value = $false
$ws.Range("J73").Value = "NA"
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = "NA"
$ws.Range("M73").Value = "NA"
$ws.Range("N73").Value = "NA"
$ws.Range("O73").Value = "NA"
$ws.Range("P73").Value = "NA"

# ---------------------------------------------------------------------
# Row 74
# ---------------------------------------------------------------------
$ws.Range("D74").Value = $false
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = $false
$ws.Range("G74").Value = 150
$ws.Range("H74").Value = $true
$ws.Range("I74").Value = $false
$ws.Range("J74").Value = "NA"
$ws.Range("K74").Value = 1
$ws.Range("L74").Value = 1
$ws.Range("M74").Value = 4
$ws.Range("N74").Value = 3
$ws.Range("O74").Value = "In terms of religious coping, those who belonged to nuclear families (16.11 [5.43] vs. 13.38 [7.17]; t-value 2.083*, [0.040]) and rural locality"
$ws.Range("P74").Value = "NA"

# ---------------------------------------------------------------------
# Row 75
# ---------------------------------------------------------------------
$ws.Range("D75").Value = $false
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = $false
$ws.Range("G75").Value = "NA"
$ws.Range("H75").Value = "NA"
$ws.Range("I75").Value = $false
$ws.Range("J75").Value = "NA"
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = "NA"
$ws.Range("M75").Value = "NA"
$ws.Range("N75").Value = "NA"
$ws.Range("O75").Value = "NA"
$ws.Range("P75").Value = "NA"

# ---------------------------------------------------------------------
# Row 76
# ---------------------------------------------------------------------
$ws.Range("D76").Value = $false
$ws.Range("E76").Value = "NA"
$ws.Range("F76").Value = $false
$ws.Range("G76").Value = "NA"
$ws.Range("H76").Value = "NA"
$ws.Range("I76").Value = $false
$ws.Range("J76").Value = "NA"
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = "NA"
$ws.Range("M76").Value = "NA"
$ws.Range("N76").Value = "NA"
$ws.Range("O76").Value = "NA"
$ws.Range("P76").Value = "NA"

# ---------------------------------------------------------------------
# Row 77 - only the sample size column is newly populated, the rest of
# the row (L77 / N77 random-number helper formulas) stay as formulas.
# ---------------------------------------------------------------------
$ws.Range("G77").Value = 27

# ---------------------------------------------------------------------
# View state: zoom + selected cell, matching where the author had
# scrolled/clicked to when saving the workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 53
$ws.Range("K77").Select()
